$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "29.928.16"
Set-TextValue "E2" "  +0.61%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.633.13"
Set-TextValue "E3" "  +1.64%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.16%  "

# Row 5 - BNB
Set-TextValue "D5" "214.46"
Set-TextValue "E5" "  +0.61%  "

# Row 6 - XRP
Set-TextValue "E6" "  +0.21%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.16%  "

# Row 8 - Solana
Set-TextValue "D8" "28.38"
Set-TextValue "E8" "  +0.49%  "

# Row 9 - Cardano
Set-TextValue "E9" "  +1.59%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  +0.73%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.866.25"
Set-TextValue "E12" "  +1.60%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.622.96"
Set-TextValue "E13" "  +0.90%  "

# Row 14 - Polygon
Set-TextValue "E14" "  +2.25%  "

# Row 15 - Chainlink
Set-TextValue "D15" "9.25"
Set-TextValue "E15" "  +17.25%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "29.934.63"
Set-TextValue "E16" "  +0.71%  "

# Row 17 - Polkadot
Set-TextValue "D17" "3.84"
Set-TextValue "E17" "  +1.98%  "

# Row 18 - Litecoin
Set-TextValue "D18" "64.01"
Set-TextValue "E18" "  -0.26%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "242.25"
Set-TextValue "E19" "  -0.16%  "

# Row 20 - ShibaInu
Set-TextValue "D20" "0.0₃0700"
Set-TextValue "E20" "  +0.17%  "

# Row 21 - Dai
Set-TextValue "E21" "  +0.07%  "

# Row 22 - Avalanche
Set-TextValue "D22" "9.86"
Set-TextValue "E22" "  +4.88%  "

# Row 23 - Uniswap
Set-TextValue "E23" "  +2.15%  "

# Row 24 - Toncoin
Set-TextValue "E24" "  +1.50%  "

# Row 25 - Monero
Set-TextValue "D25" "157.68"
Set-TextValue "E25" "  +1.43%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "15.49"
Set-TextValue "E26" "  +0.00%  "

# Row 27 - Stellar
Set-TextValue "E27" "  +0.80%  "

# Row 28 - Cosmos
Set-TextValue "E28" "  +2.06%  "

# Row 29 - BinanceUSD
Set-TextValue "E29" "  +0.10%  "

# Row 30 - Hedera
Set-TextValue "E30" "  +1.14%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.12"
Set-TextValue "E31" "  +4.15%  "

# Row 32 - Filecoin
Set-TextValue "E32" "  +3.76%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "3.17"
Set-TextValue "E33" "  -0.77%  "

# Row 34 - Maker
Set-TextValue "D34" "1.423.25"
Set-TextValue "E34" "  +0.09%  "

# Row 36 - TrustWalletToken
Set-TextValue "D36" "1.04"
Set-TextValue "E36" "  +0.04%  "

# Row 37 - MXToken
Set-TextValue "E37" "  -3.72%  "

# Row 38 - HuobiToken
Set-TextValue "E38" "  +0.03%  "

# Row 39 - VeChain
Set-TextValue "E39" "  -0.12%  "

# Row 40 - Aave
Set-TextValue "E40" "  +13.06%  "

# Row 41 - ImmutableX
Set-TextValue "E41" "  +0.87%  "

# Row 42 - RenderToken
Set-TextValue "E42" "  +1.63%  "

# Row 43 - ARBITRUM
Set-TextValue "D43" "0.827"
Set-TextValue "E43" "  +1.08%  "

# Row 44 - Kaspa
Set-TextValue "E44" "  -1.83%  "

# Row 45 - PaxDollar
Set-TextValue "E45" "  +0.13%  "

# Row 46 - now BitcoinSV (was WEMIXToken)
Set-TextValue "B46" "BitcoinSV"
Set-TextValue "C46" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D46" "52.95"
Set-TextValue "E46" "  -6.87%  "

# Row 47 - now WEMIXToken (was BitcoinSV)
Set-TextValue "B47" "WEMIXToken"
Set-TextValue "C47" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D47" "1.01"
Set-TextValue "E47" "  +3.02%  "

# Row 48 - RocketPoolETH
Set-TextValue "D48" "1.775.24"
Set-TextValue "E48" "  +1.93%  "

# Row 50 - now Quant (was BabyDogeCoin)
Set-TextValue "B50" "Quant"
Set-TextValue "C50" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D50" "89.81"
Set-TextValue "E50" "  +3.81%  "

# Row 51 - now BabyDogeCoin (was Quant)
Set-TextValue "B51" "BabyDogeCoin"
Set-TextValue "C51" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D51" "0.0₆0112"
Set-TextValue "E51" "  +9.69%  "
